$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1106.8572
$ws.Range("I28").Value = 1344.5454
$ws.Range("J28").Value = 235.33333
$ws.Range("K28").Value = 1344.5454
$ws.Range("L28").Value = 235.33333
$ws.Range("M28").Value = -859.5454
$ws.Range("N28").Value = -1205.33333

$ws.Range("H80").Value = 1005.38464
$ws.Range("I80").Value = 1314.4
$ws.Range("J80").Value = 812.25
$ws.Range("K80").Value = 3943.2
$ws.Range("L80").Value = 2436.75
$ws.Range("M80").Value = -2945.2
$ws.Range("N80").Value = -4432.75

$ws.Range("H83").Value = 1005.38464
$ws.Range("I83").Value = 1314.4
$ws.Range("J83").Value = 812.25
$ws.Range("K83").Value = 11829.6
$ws.Range("L83").Value = 7310.25
$ws.Range("M83").Value = -6837.6
$ws.Range("N83").Value = -17294.25

$ws.Range("H92").Value = 1559.5
$ws.Range("I92").Value = 1614.7693
$ws.Range("J92").Value = 1415.8
$ws.Range("K92").Value = 1614.7693
$ws.Range("L92").Value = 1415.8
$ws.Range("M92").Value = -366.7692999999999
$ws.Range("N92").Value = -3911.8

$ws.Range("H96").Value = 2073.3333
$ws.Range("I96").Value = 3372
$ws.Range("J96").Value = 450
$ws.Range("K96").Value = 10116
$ws.Range("L96").Value = 1350
$ws.Range("M96").Value = -8743
$ws.Range("N96").Value = -4096

$ws.Range("H100").Value = 1961.7222
$ws.Range("I100").Value = 1935
$ws.Range("J100").Value = 1967.0667
$ws.Range("K100").Value = 1935
$ws.Range("L100").Value = 1967.0667
$ws.Range("M100").Value = -1394
$ws.Range("N100").Value = -3049.0667

$ws.Range("H101").Value = 983.3333
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 983.3333
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 2949.9999
$ws.Range("M101").ClearContents()
$ws.Range("N101").Value = -6193.9999

$ws.Range("H132").Value = 6542221.5
$ws.Range("I132").Value = 9013320
$ws.Range("J132").Value = 11460.571
$ws.Range("K132").Value = 27039960
$ws.Range("L132").Value = 34381.713
$ws.Range("M132").Value = -27037430
$ws.Range("N132").Value = -39441.713

$ws.Range("H135").Value = 1161.6857
$ws.Range("I135").Value = 386.68967
$ws.Range("J135").Value = 4907.5
$ws.Range("K135").Value = 3480.20703
$ws.Range("L135").Value = 44167.5
$ws.Range("M135").Value = -945.2070299999996
$ws.Range("N135").Value = -49237.5

$ws.Range("H137").Value = 1288.5883
$ws.Range("I137").Value = 898.7222
$ws.Range("J137").Value = 1727.1875
$ws.Range("K137").Value = 2696.1666
$ws.Range("L137").Value = 5181.5625
$ws.Range("M137").Value = -146.1666
$ws.Range("N137").Value = -10281.5625

$ws.Range("H138").Value = 1481.404
$ws.Range("I138").Value = 762.3182
$ws.Range("J138").Value = 1686.8572
$ws.Range("K138").Value = 2286.9546
$ws.Range("L138").Value = 5060.571599999999
$ws.Range("M138").Value = 2853.0454
$ws.Range("N138").Value = -15340.5716

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4028.3635
$ws.Range("I32").Value = 3522.65
$ws.Range("J32").Value = 9085.5
$ws.Range("K32").Value = 3522.65
$ws.Range("L32").Value = 9085.5
$ws.Range("M32").Value = -3235.65
$ws.Range("N32").Value = -9659.5

$ws.Range("H55").Value = 37500
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 37500
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 37500
$ws.Range("M55").ClearContents()
$ws.Range("N55").Value = -38130

$ws.Range("H74").Value = 1573.3611
$ws.Range("I74").Value = 854.75
$ws.Range("J74").Value = 2471.625
$ws.Range("K74").Value = 854.75
$ws.Range("L74").Value = 2471.625
$ws.Range("M74").Value = 19.25
$ws.Range("N74").Value = -4219.625

$ws.Range("H77").Value = 1573.3611
$ws.Range("I77").Value = 854.75
$ws.Range("J77").Value = 2471.625
$ws.Range("K77").Value = 4273.75
$ws.Range("L77").Value = 12358.125
$ws.Range("M77").Value = 94.25
$ws.Range("N77").Value = -21094.125

$ws.Range("H97").Value = 682.8
$ws.Range("I97").Value = 404.66666
$ws.Range("J97").Value = 1100
$ws.Range("K97").Value = 404.66666
$ws.Range("L97").Value = 1100
$ws.Range("M97").Value = 91.33334000000002
$ws.Range("N97").Value = -2092

$ws.Range("H122").Value = 1040.2
$ws.Range("I122").Value = 863.05
$ws.Range("J122").Value = 1748.8
$ws.Range("K122").Value = 2589.15
$ws.Range("L122").Value = 5246.4
$ws.Range("M122").Value = -139.1499999999996
$ws.Range("N122").Value = -10146.4

$ws.Range("H132").Value = 1537.2826
$ws.Range("I132").Value = 1253.5161
$ws.Range("J132").Value = 2123.7334
$ws.Range("K132").Value = 3760.5483
$ws.Range("L132").Value = 6371.2002
$ws.Range("M132").Value = -1230.5483
$ws.Range("N132").Value = -11431.2002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1702.2858
$ws.Range("I31").Value = 1470.2858
$ws.Range("J31").Value = 1934.2858
$ws.Range("K31").Value = 1470.2858
$ws.Range("L31").Value = 1934.2858
$ws.Range("M31").Value = -1175.2858
$ws.Range("N31").Value = -2524.2858

$ws.Range("H34").Value = 1702.2858
$ws.Range("I34").Value = 1470.2858
$ws.Range("J34").Value = 1934.2858
$ws.Range("K34").Value = 1470.2858
$ws.Range("L34").Value = 1934.2858
$ws.Range("M34").Value = -1268.2858
$ws.Range("N34").Value = -2338.2858

$ws.Range("H58").Value = 1000.8947
$ws.Range("I58").Value = 898.4375
$ws.Range("J58").Value = 1547.3334
$ws.Range("K58").Value = 898.4375
$ws.Range("L58").Value = 1547.3334
$ws.Range("M58").Value = -695.4375
$ws.Range("N58").Value = -1953.3334

$ws.Range("H132").Value = 1377.4807
$ws.Range("I132").Value = 1013.7174
$ws.Range("J132").Value = 4166.3335
$ws.Range("K132").Value = 3041.1522
$ws.Range("L132").Value = 12499.0005
$ws.Range("M132").Value = -511.1522
$ws.Range("N132").Value = -17559.0005

$ws.Range("H136").Value = 1000.8947
$ws.Range("I136").Value = 898.4375
$ws.Range("J136").Value = 1547.3334
$ws.Range("K136").Value = 2695.3125
$ws.Range("L136").Value = 4642.0002
$ws.Range("M136").Value = -145.3125
$ws.Range("N136").Value = -9742.0002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 8392.904
$ws.Range("I3").Value = 3873.8462
$ws.Range("J3").Value = 15736.375
$ws.Range("K3").Value = 11621.5386
$ws.Range("L3").Value = 47209.125
$ws.Range("M3").Value = -11509.5386
$ws.Range("N3").Value = -47433.125

$ws.Range("H33").Value = 303.33334
$ws.Range("I33").Value = 214.54546
$ws.Range("J33").Value = 442.85715
$ws.Range("K33").Value = 1287.27276
$ws.Range("L33").Value = 2657.1429
$ws.Range("M33").Value = -1004.27276
$ws.Range("N33").Value = -3223.1429

$ws.Range("H44").Value = 2920.8
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 2920.8
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 8762.400000000001
$ws.Range("M44").ClearContents()
$ws.Range("N44").Value = -9558.400000000001

$ws.Range("H107").Value = 6285.647
$ws.Range("I107").Value = 391.4
$ws.Range("J107").Value = 8741.583000000001
$ws.Range("K107").Value = 1174.2
$ws.Range("L107").Value = 26224.749
$ws.Range("M107").Value = 745.8000000000002
$ws.Range("N107").Value = -30064.749

$ws.Range("H125").Value = 2676.6667
$ws.Range("I125").Value = 998.5
$ws.Range("J125").Value = 6033
$ws.Range("K125").Value = 2995.5
$ws.Range("L125").Value = 18099
$ws.Range("M125").Value = 1924.5
$ws.Range("N125").Value = -27939

$ws.Range("H127").Value = 10000
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 10000
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 30000
$ws.Range("M127").ClearContents()
$ws.Range("N127").Value = -39920

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H41").Value = 1333.3334
$ws.Range("I41").Value = 2000
$ws.Range("J41").Value = 1000
$ws.Range("K41").Value = 2000
$ws.Range("L41").Value = 1000
$ws.Range("M41").Value = -1645
$ws.Range("N41").Value = -1710

$ws.Range("H57").Value = 22250
$ws.Range("I57").Value = 22000
$ws.Range("J57").Value = 22333.334
$ws.Range("K57").Value = 22000
$ws.Range("L57").Value = 22333.334
$ws.Range("M57").Value = -21180
$ws.Range("N57").Value = -23973.334

$ws.Range("H70").Value = 64302700
$ws.Range("I70").Value = 62527348
$ws.Range("J70").Value = 66669830
$ws.Range("K70").Value = 62527348
$ws.Range("L70").Value = 66669830
$ws.Range("M70").Value = -62527078
$ws.Range("N70").Value = -66670370

$ws.Range("H73").Value = 64302700
$ws.Range("I73").Value = 62527348
$ws.Range("J73").Value = 66669830
$ws.Range("K73").Value = 62527348
$ws.Range("L73").Value = 66669830
$ws.Range("M73").Value = -62526412
$ws.Range("N73").Value = -66671702

$ws.Range("H102").Value = 5238.5713
$ws.Range("I102").Value = 3759.9473
$ws.Range("J102").Value = 8360.111000000001
$ws.Range("K102").Value = 3759.9473
$ws.Range("L102").Value = 8360.111000000001
$ws.Range("M102").Value = -2137.9473
$ws.Range("N102").Value = -11604.111

$ws.Range("H113").Value = 2000
$ws.Range("I113").Value = 2000
$ws.Range("J113").Value = 0
$ws.Range("K113").ClearContents()
$ws.Range("L113").ClearContents()
$ws.Range("M113").Value = 170
$ws.Range("N113").Value = 0

$ws.Range("H126").Value = 2220
$ws.Range("I126").Value = 1908
$ws.Range("J126").Value = 3000
$ws.Range("K126").Value = 5724
$ws.Range("L126").Value = 9000
$ws.Range("M126").Value = -3254
$ws.Range("N126").Value = -13940

$ws.Range("H132").Value = 2091
$ws.Range("I132").Value = 1342.8948
$ws.Range("J132").Value = 3383.182
$ws.Range("K132").Value = 4028.6844
$ws.Range("L132").Value = 10149.546
$ws.Range("M132").Value = -1498.6844
$ws.Range("N132").Value = -15209.546

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 865.7692
$ws.Range("I22").Value = 759.1667
$ws.Range("J22").Value = 957.1429000000001
$ws.Range("K22").Value = 759.1667
$ws.Range("L22").Value = 957.1429000000001
$ws.Range("M22").Value = -464.1667
$ws.Range("N22").Value = -1547.1429

$ws.Range("H27").Value = 865.7692
$ws.Range("I27").Value = 759.1667
$ws.Range("J27").Value = 957.1429000000001
$ws.Range("K27").Value = 759.1667
$ws.Range("L27").Value = 957.1429000000001
$ws.Range("M27").Value = -652.1667
$ws.Range("N27").Value = -1171.1429

$ws.Range("H55").Value = 384.9375
$ws.Range("I55").Value = 286.4
$ws.Range("J55").Value = 549.1667
$ws.Range("K55").Value = 286.4
$ws.Range("L55").Value = 549.1667
$ws.Range("M55").Value = -113.4
$ws.Range("N55").Value = -895.1667

$ws.Range("H93").Value = 784.8461
$ws.Range("I93").Value = 745.36365
$ws.Range("J93").Value = 1002
$ws.Range("K93").Value = 745.36365
$ws.Range("L93").Value = 1002
$ws.Range("M93").Value = 502.63635
$ws.Range("N93").Value = -3498

$ws.Range("H122").Value = 27780534
$ws.Range("I122").Value = 62501950
$ws.Range("J122").Value = 3400
$ws.Range("K122").Value = 187505850
$ws.Range("L122").Value = 10200
$ws.Range("M122").Value = -187503400
$ws.Range("N122").Value = -15100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H49").Value = 5000
$ws.Range("I49").Value = 5000
$ws.Range("J49").Value = 5000
$ws.Range("K49").Value = 5000
$ws.Range("L49").Value = 5000
$ws.Range("M49").Value = -4770
$ws.Range("N49").Value = -5460

$ws.Range("H103").Value = 12000
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 12000
$ws.Range("K103").Value = 0
$ws.Range("L103").Value = 12000
$ws.Range("M103").ClearContents()
$ws.Range("N103").Value = -14344

$ws.Range("H104").Value = 19980
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 19980
$ws.Range("K104").Value = 0
$ws.Range("L104").Value = 19980
$ws.Range("M104").ClearContents()
$ws.Range("N104").Value = -26968

$ws.Range("H105").Value = 36000
$ws.Range("I105").Value = 36000
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 36000
$ws.Range("L105").ClearContents()
$ws.Range("M105").Value = -32506
$ws.Range("N105").Value = 0
